# Generate Report for Handback
#
# This script mirrors the "handback" run of the localization-status report:
#   - the Status column flips from "In Translation" to
#     "Handed back: in sync with en-US" for every tracked file,
#   - each language sheet (zh-cn / de-de) now has a resolved Target File +
#     Handback File + Handback DateTime for both tracked source files,
#   - those two new file-name cells (column I) become hyperlinks, just like
#     column A already is,
#   - a handful of columns are widened so the newly-populated file-name /
#     status text is readable.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$file1Name = "143163b4-64bd-477b-86ca-ed8e56ef20e6.md"
$file1Url  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4da697efc48ec40db91167b45982c03055c64c7/e2e/143163b4-64bd-477b-86ca-ed8e56ef20e6.md"
$file1ZhHandback = "143163b4-64bd-477b-86ca-ed8e56ef20e6.99d9cf5e04f217c55e20e2970a3696f4bd5daaf0.zh-cn.xlf"
$file1DeHandback = "143163b4-64bd-477b-86ca-ed8e56ef20e6.99d9cf5e04f217c55e20e2970a3696f4bd5daaf0.de-de.xlf"

$file2Name = "e58989c1-4165-4e19-927b-997daff9e821.md"
$file2Url  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4da697efc48ec40db91167b45982c03055c64c7/e2e/e58989c1-4165-4e19-927b-997daff9e821.md"
$file2ZhHandback = "e58989c1-4165-4e19-927b-997daff9e821.2fe87d0478eefe739d64a4c0a26df83e84adc2b6.zh-cn.xlf"
$file2DeHandback = "e58989c1-4165-4e19-927b-997daff9e821.2fe87d0478eefe739d64a4c0a26df83e84adc2b6.de-de.xlf"

$zhHandbackDate = "2016-09-04 18:28:31"
$deHandbackDate = "2016-09-04 18:28:39"

# Cornflower blue (FF6495ED), matches the workbook's existing "HyperLink" style.
$hyperlinkColor = 15570276

function Set-HyperlinkLook($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Overview sheet: its zh-cn / de-de status cells share the same
# "In Translation" shared string, so they flip to the new text too.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("I2").Value = $file1Name
$wsZh.Range("J2").Value = $file1ZhHandback
$wsZh.Range("K2").Value = $zhHandbackDate

$wsZh.Range("I3").Value = $file2Name
$wsZh.Range("J3").Value = $file2ZhHandback
$wsZh.Range("K3").Value = $zhHandbackDate

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $file1Url, "", "", $file1Name)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $file2Url, "", "", $file2Name)

Set-HyperlinkLook $wsZh.Range("I2")
Set-HyperlinkLook $wsZh.Range("I3")

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("I2").Value = $file1Name
$wsDe.Range("J2").Value = $file1DeHandback
$wsDe.Range("K2").Value = $deHandbackDate

$wsDe.Range("I3").Value = $file2Name
$wsDe.Range("J3").Value = $file2DeHandback
$wsDe.Range("K3").Value = $deHandbackDate

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $file1Url, "", "", $file1Name)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $file2Url, "", "", $file2Name)

Set-HyperlinkLook $wsDe.Range("I2")
Set-HyperlinkLook $wsDe.Range("I3")

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664
